# Scheduled refresh of market-board price data across the per-job
# "Ragnarok_Profits" sheets. Updates currentAveragePrice / NQ / HQ price
# columns (H:N) for the leves whose listings moved since the last run.
# Column layout per table: H=currentAveragePrice, I=currentAveragePriceNQ,
# J=currentAveragePriceHQ, K=LevePriceNQ, L=LevePriceHQ, M=LeveProfitNQ,
# N=LeveProfitHQ.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 71430780
$ws.Cells.Item(40, 10).Value = 250001980
$ws.Cells.Item(40, 12).Value = 250001980
$ws.Cells.Item(40, 14).Value = -250002330

$ws.Cells.Item(70, 8).Value = 92018780
$ws.Cells.Item(70, 9).Value = 202439630
$ws.Cells.Item(70, 10).Value = 1415
$ws.Cells.Item(70, 11).Value = 607318890
$ws.Cells.Item(70, 12).Value = 4245
$ws.Cells.Item(70, 13).Value = -607318620
$ws.Cells.Item(70, 14).Value = -4785

$ws.Cells.Item(73, 8).Value = 92018780
$ws.Cells.Item(73, 9).Value = 202439630
$ws.Cells.Item(73, 10).Value = 1415
$ws.Cells.Item(73, 11).Value = 607318890
$ws.Cells.Item(73, 12).Value = 4245
$ws.Cells.Item(73, 13).Value = -607317954
$ws.Cells.Item(73, 14).Value = -6117

$ws.Cells.Item(99, 8).Value = 901.375
$ws.Cells.Item(99, 10).Value = 1141.1666
$ws.Cells.Item(99, 12).Value = 3423.4998
$ws.Cells.Item(99, 14).Value = -6419.4998

$ws.Cells.Item(111, 8).Value = 11555.714
$ws.Cells.Item(111, 9).Value = 4328
$ws.Cells.Item(111, 10).Value = 16976.5
$ws.Cells.Item(111, 11).Value = 12984
$ws.Cells.Item(111, 12).Value = 50929.5
$ws.Cells.Item(111, 13).Value = -9917
$ws.Cells.Item(111, 14).Value = -57063.5

$ws.Cells.Item(123, 8).Value = 99998.5
$ws.Cells.Item(123, 10).Value = 99998.5
$ws.Cells.Item(123, 12).Value = 99998.5
$ws.Cells.Item(123, 14).Value = -109798.5

$ws.Cells.Item(132, 8).Value = 3825.9111
$ws.Cells.Item(132, 9).Value = 2671.262
$ws.Cells.Item(132, 10).Value = 19991
$ws.Cells.Item(132, 11).Value = 8013.786
$ws.Cells.Item(132, 12).Value = 59973
$ws.Cells.Item(132, 13).Value = -5483.786
$ws.Cells.Item(132, 14).Value = -65033

$ws.Cells.Item(135, 8).Value = 1181.5758
$ws.Cells.Item(135, 9).Value = 1037.5416
$ws.Cells.Item(135, 10).Value = 1565.6666
$ws.Cells.Item(135, 11).Value = 9337.874400000001
$ws.Cells.Item(135, 12).Value = 14090.9994
$ws.Cells.Item(135, 13).Value = -6802.874400000001
$ws.Cells.Item(135, 14).Value = -19160.9994

$ws.Cells.Item(138, 8).Value = 7753.125
$ws.Cells.Item(138, 9).Value = 11696
$ws.Cells.Item(138, 10).Value = 7189.857
$ws.Cells.Item(138, 11).Value = 35088
$ws.Cells.Item(138, 12).Value = 21569.571
$ws.Cells.Item(138, 13).Value = -29948
$ws.Cells.Item(138, 14).Value = -31849.571

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 2664.155
$ws.Cells.Item(32, 9).Value = 2799.9846
$ws.Cells.Item(32, 11).Value = 2799.9846
$ws.Cells.Item(32, 13).Value = -2512.9846

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 1806.75
$ws.Cells.Item(16, 9).Value = 1697.3846
$ws.Cells.Item(16, 10).Value = 2280.6667
$ws.Cells.Item(16, 11).Value = 1697.3846
$ws.Cells.Item(16, 12).Value = 2280.6667
$ws.Cells.Item(16, 13).Value = -1410.3846
$ws.Cells.Item(16, 14).Value = -2854.6667

$ws.Cells.Item(31, 8).Value = 2695
$ws.Cells.Item(31, 9).Value = 2617.7368
$ws.Cells.Item(31, 10).Value = 2807.923
$ws.Cells.Item(31, 11).Value = 2617.7368
$ws.Cells.Item(31, 12).Value = 2807.923
$ws.Cells.Item(31, 13).Value = -2322.7368
$ws.Cells.Item(31, 14).Value = -3397.923

$ws.Cells.Item(34, 8).Value = 2695
$ws.Cells.Item(34, 9).Value = 2617.7368
$ws.Cells.Item(34, 10).Value = 2807.923
$ws.Cells.Item(34, 11).Value = 2617.7368
$ws.Cells.Item(34, 12).Value = 2807.923
$ws.Cells.Item(34, 13).Value = -2415.7368
$ws.Cells.Item(34, 14).Value = -3211.923

$ws.Cells.Item(105, 8).Value = 1294.5834
$ws.Cells.Item(105, 9).Value = 1066.875
$ws.Cells.Item(105, 10).Value = 1750
$ws.Cells.Item(105, 11).Value = 1066.875
$ws.Cells.Item(105, 12).Value = 1750
$ws.Cells.Item(105, 13).Value = 680.125
$ws.Cells.Item(105, 14).Value = -5244

$ws.Cells.Item(113, 8).Value = 1806.75
$ws.Cells.Item(113, 9).Value = 1697.3846
$ws.Cells.Item(113, 10).Value = 2280.6667
$ws.Cells.Item(113, 11).Value = 1697.3846
$ws.Cells.Item(113, 12).Value = 2280.6667
$ws.Cells.Item(113, 13).Value = 472.6153999999999
$ws.Cells.Item(113, 14).Value = -6620.6667

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(26, 8).Value = 5160.4287
$ws.Cells.Item(26, 9).Value = 538
$ws.Cells.Item(26, 10).Value = 16716.5
$ws.Cells.Item(26, 11).Value = 1614
$ws.Cells.Item(26, 12).Value = 50149.5
$ws.Cells.Item(26, 13).Value = -1326
$ws.Cells.Item(26, 14).Value = -50725.5

$ws.Cells.Item(98, 8).Value = 698
$ws.Cells.Item(98, 10).Value = 677
$ws.Cells.Item(98, 12).Value = 2031
$ws.Cells.Item(98, 14).Value = -5027

$ws.Cells.Item(103, 8).Value = 4160.9165
$ws.Cells.Item(103, 9).Value = 333.22223
$ws.Cells.Item(103, 10).Value = 15644
$ws.Cells.Item(103, 11).Value = 999.66669
$ws.Cells.Item(103, 12).Value = 46932
$ws.Cells.Item(103, 13).Value = -120.66669
$ws.Cells.Item(103, 14).Value = -48690

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(15, 8).Value = 46332.668
$ws.Cells.Item(15, 10).Value = 46332.668
$ws.Cells.Item(15, 12).Value = 46332.668
$ws.Cells.Item(15, 14).Value = -46908.668

$ws.Cells.Item(81, 8).Value = 46332.668
$ws.Cells.Item(81, 10).Value = 46332.668
$ws.Cells.Item(81, 12).Value = 46332.668
$ws.Cells.Item(81, 14).Value = -48328.668

$ws.Cells.Item(84, 8).Value = 46332.668
$ws.Cells.Item(84, 10).Value = 46332.668
$ws.Cells.Item(84, 12).Value = 138998.004
$ws.Cells.Item(84, 14).Value = -148982.004

$ws.Cells.Item(122, 8).Value = 3296.2856
$ws.Cells.Item(122, 9).Value = 3679
$ws.Cells.Item(122, 10).Value = 1000
$ws.Cells.Item(122, 11).Value = 11037
$ws.Cells.Item(122, 12).Value = 3000
$ws.Cells.Item(122, 13).Value = -8587
$ws.Cells.Item(122, 14).Value = -7900

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 1274.5
$ws.Cells.Item(46, 9).Value = 1300
$ws.Cells.Item(46, 10).Value = 1249
$ws.Cells.Item(46, 11).Value = 1300
$ws.Cells.Item(46, 12).Value = 1249
$ws.Cells.Item(46, 13).Value = -1112
$ws.Cells.Item(46, 14).Value = -1625

$ws.Cells.Item(55, 8).Value = 1322.0555
$ws.Cells.Item(55, 9).Value = 1113.8572
$ws.Cells.Item(55, 10).Value = 1454.5454
$ws.Cells.Item(55, 11).Value = 1113.8572
$ws.Cells.Item(55, 12).Value = 1454.5454
$ws.Cells.Item(55, 13).Value = -940.8571999999999
$ws.Cells.Item(55, 14).Value = -1800.5454

$ws.Cells.Item(68, 8).Value = 13891157
$ws.Cells.Item(68, 9).Value = 41666664
$ws.Cells.Item(68, 11).Value = 41666664
$ws.Cells.Item(68, 13).Value = -41665915

$ws.Cells.Item(71, 8).Value = 13891157
$ws.Cells.Item(71, 9).Value = 41666664
$ws.Cells.Item(71, 11).Value = 208333320
$ws.Cells.Item(71, 13).Value = -208329576

$ws.Cells.Item(100, 8).Value = 12517897
$ws.Cells.Item(100, 9).Value = 6169.5
$ws.Cells.Item(100, 10).Value = 25029624
$ws.Cells.Item(100, 11).Value = 6169.5
$ws.Cells.Item(100, 12).Value = 25029624
$ws.Cells.Item(100, 13).Value = -5628.5
$ws.Cells.Item(100, 14).Value = -25030706

$ws.Cells.Item(132, 8).Value = 2019.3823
$ws.Cells.Item(132, 9).Value = 1958.4062
$ws.Cells.Item(132, 11).Value = 5875.2186
$ws.Cells.Item(132, 13).Value = -3345.2186

$ws.Cells.Item(136, 8).Value = 2408.276
$ws.Cells.Item(136, 9).Value = 2348.1155
$ws.Cells.Item(136, 11).Value = 7044.3465
$ws.Cells.Item(136, 13).Value = -4494.3465

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 2575
$ws.Cells.Item(122, 9).Value = 2575
$ws.Cells.Item(122, 10).Value = 0
$ws.Cells.Item(122, 11).Value = 7725
$ws.Cells.Item(122, 12).Value = 0
$ws.Cells.Item(122, 13).Value = -5275
$ws.Cells.Item(122, 14).ClearContents()

$ws.Cells.Item(136, 8).Value = 239515.77
$ws.Cells.Item(136, 10).Value = 1434733.6
$ws.Cells.Item(136, 12).Value = 4304200.800000001
$ws.Cells.Item(136, 14).Value = -4309300.800000001
